# gsc-export/HTTPS.xlsx -- updated main GSC export data
# The daily export window rolls forward by one day: the oldest date
# (2025-10-25) drops off the top, every remaining date/count shifts up
# one row, and a new row for 2026-01-23 is appended at the bottom.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# New values for Date (col A) and HTTPS "Pages" (col C), rows 2-91.
$dates = @(
    "2025-10-26",
    "2025-10-27",
    "2025-10-28",
    "2025-10-29",
    "2025-10-30",
    "2025-10-31",
    "2025-11-01",
    "2025-11-02",
    "2025-11-03",
    "2025-11-04",
    "2025-11-05",
    "2025-11-06",
    "2025-11-07",
    "2025-11-08",
    "2025-11-09",
    "2025-11-10",
    "2025-11-11",
    "2025-11-12",
    "2025-11-13",
    "2025-11-14",
    "2025-11-15",
    "2025-11-16",
    "2025-11-17",
    "2025-11-18",
    "2025-11-19",
    "2025-11-20",
    "2025-11-21",
    "2025-11-22",
    "2025-11-23",
    "2025-11-24",
    "2025-11-25",
    "2025-11-26",
    "2025-11-27",
    "2025-11-28",
    "2025-11-29",
    "2025-11-30",
    "2025-12-01",
    "2025-12-02",
    "2025-12-03",
    "2025-12-04",
    "2025-12-05",
    "2025-12-06",
    "2025-12-07",
    "2025-12-08",
    "2025-12-09",
    "2025-12-10",
    "2025-12-11",
    "2025-12-12",
    "2025-12-13",
    "2025-12-14",
    "2025-12-15",
    "2025-12-16",
    "2025-12-17",
    "2025-12-18",
    "2025-12-19",
    "2025-12-20",
    "2025-12-21",
    "2025-12-22",
    "2025-12-23",
    "2025-12-24",
    "2025-12-25",
    "2025-12-26",
    "2025-12-27",
    "2025-12-28",
    "2025-12-29",
    "2025-12-30",
    "2025-12-31",
    "2026-01-01",
    "2026-01-02",
    "2026-01-03",
    "2026-01-04",
    "2026-01-05",
    "2026-01-06",
    "2026-01-07",
    "2026-01-08",
    "2026-01-09",
    "2026-01-10",
    "2026-01-11",
    "2026-01-12",
    "2026-01-13",
    "2026-01-14",
    "2026-01-15",
    "2026-01-16",
    "2026-01-17",
    "2026-01-18",
    "2026-01-19",
    "2026-01-20",
    "2026-01-21",
    "2026-01-22",
    "2026-01-23"
)

$pages = @(
    86,
    90,
    83,
    90,
    93,
    92,
    102,
    113,
    115,
    107,
    105,
    100,
    94,
    86,
    83,
    66,
    54,
    46,
    43,
    40,
    37,
    35,
    30,
    29,
    26,
    25,
    25,
    26,
    26,
    25,
    25,
    27,
    28,
    28,
    27,
    27,
    27,
    27,
    27,
    26,
    25,
    25,
    25,
    26,
    27,
    27,
    29,
    29,
    30,
    30,
    31,
    31,
    31,
    31,
    31,
    32,
    32,
    32,
    32,
    30,
    31,
    32,
    30,
    28,
    28,
    28,
    28,
    29,
    29,
    28,
    27,
    27,
    28,
    27,
    27,
    27,
    27,
    26,
    26,
    27,
    26,
    26,
    25,
    25,
    25,
    25,
    26,
    25,
    24,
    23
)

# Column A holds the dates as plain text (shared strings), not Excel date
# serials. Format the range as Text first so assigning date-look-alike
# strings does not get auto-converted to a date value, then clear that
# temporary formatting override once all the values are written so the
# cells keep the workbook default style.
$dateRange = $ws.Range("A2:A91")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 3).Value = $pages[$i]
}

$dateRange.ClearFormats()
